# Append the latest Nalco PDF run-log entry (row 50) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 50

# Match the existing data-row formatting (center/center alignment, same as
# every other logged run in A2:H49).
$dataRange = $ws.Range("A$row`:H$row")
$dataRange.HorizontalAlignment = -4108   # xlCenter
$dataRange.VerticalAlignment = -4108     # xlCenter

$ws.Range("A$row").Value = "2025-08-24 03:58:52 UTC"
$ws.Range("B$row").Value = "2025-08-24 09:28:52 IST"
$ws.Range("C$row").Value = "SKIPPED"
$ws.Range("D$row").Value = "No change in PDF. Skipping download & Excel update."
$ws.Range("E$row").Value = "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf"
$ws.Range("F$row").Value = ""
$ws.Range("G$row").Value = 0
$ws.Range("H$row").Value = ""
